$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-12-18"

# Update the header text (shared string) for the 2022 column
$ws.Range("I1").Value = "2022 (through 12-18)"

# Update December (row 13) and Total (row 14) values for the 2022 column (I)
$ws.Range("I13").Value = 78
$ws.Range("I14").Value = 1595
